$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting existing rows 25:144 down to 26:145
$ws.Rows(25).Insert()

# Populate the new row 25 with the new record's values
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 44701
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 100112012
$ws.Range("G25").Value = "Espinaca"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 90
$ws.Range("K25").Value = 8000
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = 8556
$ws.Range("N25").Value = '$/docena de atados'
$ws.Range("O25").Value = "Región de La Araucanía"
$ws.Range("P25").Value = 2852
$ws.Range("Q25").Value = 3
$ws.Range("R25").Value = "Hortaliza"
